# Apply update to chla standard curve summary:
# 1. Add new column K "slope_rfu_over_ugL" header
# 2. Move the old slope value (F column) for "straight" rows into K,
#    and replace F with a formula =1/K<row>
# 3. For rows where F was not modified to a formula, simply copy the F value into K

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in K1
$ws.Range("K1").Value = "slope_rfu_over_ugL"

# Rows 2,3,4,7,8,9: move existing F value to K, then set F to formula 1/K
$rowsWithFormula = 2,3,4,7,8,9
foreach ($r in $rowsWithFormula) {
    $oldVal = $ws.Cells.Item($r, 6).Value2   # column F = 6
    $ws.Cells.Item($r, 11).Value = $oldVal   # column K = 11
    $ws.Range("F$r").Formula = "=1/K$r"
}

# Rows 5,6,10,11: copy F value as literal into K (no formula change to F)
$rowsLiteralCopy = 5,6,10,11
foreach ($r in $rowsLiteralCopy) {
    $oldVal = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 11).Value = $oldVal
}

# Update selection to reflect the final state (F11 selected)
$ws.Range("F11").Select()
